# Update "want to go" counts (column F) on the "展览" (Exhibitions) and
# "全部类型" (All types) sheets to reflect the latest scraped numbers.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 2-5 correspond to F3 in the same order as "全部类型"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 272
$wsExhibit.Range("F3").Value = 90
$wsExhibit.Range("F4").Value = 1006
$wsExhibit.Range("F5").Value = 555

# Sheet "全部类型" - has an extra row (卡农...) inserted at row 5, so the
# matching "展览" row 5 value lands on row 6 here.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 272
$wsAll.Range("F3").Value = 90
$wsAll.Range("F4").Value = 1006
$wsAll.Range("F6").Value = 555
